$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$normalStyle = $ws.Range("A1").Style

# D2, E2, F2 look like dates and L2 looks like a plain number -- force
# them to stay text (matching the source data) by flipping to a text
# number format before writing, then restoring the original (default)
# cell style so the published style table is unaffected.
$ws.Range("D2:F2").NumberFormat = "@"
$ws.Range("L2").NumberFormat = "@"

$ws.Range("A2").Value = "abcde"
$ws.Range("B2").Value = "First camp"
$ws.Range("C2").Value = $true
$ws.Range("D2").Value = "2023-11-26"
$ws.Range("E2").Value = "2023-01-02"
$ws.Range("F2").Value = "2023-11-29"
$ws.Range("G2").Value = "Nanyang"
$ws.Range("H2").Value = 49
$ws.Range("I2").Value = 10
$ws.Range("J2").Value = "Halloween Camp"
$ws.Range("K2").Value = "98d0e59407f946b7aed49150ceba8627"
$ws.Range("L2").Value = "1"
$ws.Range("M2").Value = $true

$ws.Range("D2:F2").Style = $normalStyle
$ws.Range("L2").Style = $normalStyle
